$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 545
$ws.Range("I31").Value = 545
$ws.Range("K31").Value = 1635
$ws.Range("M31").Value = -1405
$ws.Range("H33").Value = 599.5
$ws.Range("I33").Value = 199
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 199
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = 30
$ws.Range("N33").Value = -1458
$ws.Range("H43").Value = 6140
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10138
$ws.Range("H92").Value = 307.8889
$ws.Range("J92").Value = 118.25
$ws.Range("L92").Value = 118.25
$ws.Range("N92").Value = -2614.25
$ws.Range("H132").Value = 1349.8837
$ws.Range("I132").Value = 1201.15
$ws.Range("K132").Value = 3603.45
$ws.Range("M132").Value = -1073.45
$ws.Range("H138").Value = 34484412
$ws.Range("J138").Value = 52633384
$ws.Range("L138").Value = 157900152
$ws.Range("N138").Value = -157910432
$ws.Range("H141").Value = 3665.923
$ws.Range("I141").Value = 3314.2727
$ws.Range("J141").Value = 5600
$ws.Range("K141").Value = 9942.8181
$ws.Range("L141").Value = 16800
$ws.Range("M141").Value = -4762.8181
$ws.Range("N141").Value = -27160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 6697283.5
$ws.Range("J11").Value = 27900
$ws.Range("L11").Value = 27900
$ws.Range("N11").Value = -28188
$ws.Range("H32").Value = 9031.893
$ws.Range("I32").Value = 4891.2324
$ws.Range("K32").Value = 4891.2324
$ws.Range("M32").Value = -4604.2324
$ws.Range("H61").Value = 58072.777
$ws.Range("I61").Value = 1579.4445
$ws.Range("K61").Value = 1579.4445
$ws.Range("M61").Value = -1367.4445
$ws.Range("H74").Value = 26238.195
$ws.Range("I74").Value = 35800.516
$ws.Range("J74").Value = 3129.25
$ws.Range("K74").Value = 35800.516
$ws.Range("L74").Value = 3129.25
$ws.Range("M74").Value = -34926.516
$ws.Range("N74").Value = -4877.25
$ws.Range("H77").Value = 26238.195
$ws.Range("I77").Value = 35800.516
$ws.Range("J77").Value = 3129.25
$ws.Range("K77").Value = 179002.58
$ws.Range("L77").Value = 15646.25
$ws.Range("M77").Value = -174634.58
$ws.Range("N77").Value = -24382.25
$ws.Range("H102").Value = 66825.06
$ws.Range("I102").Value = 93332
$ws.Range("K102").Value = 93332
$ws.Range("M102").Value = -91710
$ws.Range("H122").Value = 4509.4736
$ws.Range("I122").Value = 5197.4287
$ws.Range("K122").Value = 15592.2861
$ws.Range("M122").Value = -13142.2861
$ws.Range("H135").Value = 85842.28999999999
$ws.Range("J135").Value = 85842.28999999999
$ws.Range("L135").Value = 85842.28999999999
$ws.Range("N135").Value = -95982.28999999999
$ws.Range("H136").Value = 58072.777
$ws.Range("I136").Value = 1579.4445
$ws.Range("K136").Value = 4738.333500000001
$ws.Range("M136").Value = -2188.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 62011.844
$ws.Range("I20").Value = 89194.38
$ws.Range("K20").Value = 89194.38
$ws.Range("M20").Value = -88947.38
$ws.Range("H86").Value = 5831.3335
$ws.Range("I86").Value = 5831.3335
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5831.3335
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -4708.3335
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 5831.3335
$ws.Range("I89").Value = 5831.3335
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 29156.6675
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -23540.6675
$ws.Range("N89").ClearContents()
$ws.Range("H122").Value = 74264.3
$ws.Range("J122").Value = 74264.3
$ws.Range("L122").Value = 74264.3
$ws.Range("N122").Value = -84064.3
$ws.Range("H134").Value = 1815.4286
$ws.Range("I134").Value = 1311.862
$ws.Range("J134").Value = 4249.3335
$ws.Range("K134").Value = 3935.586
$ws.Range("L134").Value = 12748.0005
$ws.Range("M134").Value = -1400.586
$ws.Range("N134").Value = -17818.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2312.1428
$ws.Range("I58").Value = 1592.2
$ws.Range("K58").Value = 1592.2
$ws.Range("M58").Value = -1389.2
$ws.Range("H107").Value = 659.4
$ws.Range("I107").Value = 730.9231
$ws.Range("K107").Value = 730.9231
$ws.Range("M107").Value = 1189.0769
$ws.Range("H134").Value = 51942.7
$ws.Range("I134").Value = 1821.3125
$ws.Range("K134").Value = 5463.9375
$ws.Range("M134").Value = -2928.9375
$ws.Range("H136").Value = 2312.1428
$ws.Range("I136").Value = 1592.2
$ws.Range("K136").Value = 4776.6
$ws.Range("M136").Value = -2226.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 262.5
$ws.Range("I25").Value = 262.5
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 787.5
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -618.5
$ws.Range("N25").ClearContents()
$ws.Range("H30").Value = 262.5
$ws.Range("I30").Value = 262.5
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 787.5
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -685.5
$ws.Range("N30").ClearContents()
$ws.Range("H131").Value = 30834.646
$ws.Range("I131").Value = 83775.664
$ws.Range("J131").Value = 1957.7273
$ws.Range("K131").Value = 251326.992
$ws.Range("L131").Value = 5873.1819
$ws.Range("M131").Value = -246286.992
$ws.Range("N131").Value = -15953.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6479826
$ws.Range("I11").Value = 7052831
$ws.Range("K11").Value = 7052831
$ws.Range("M11").Value = -7052692
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H42").Value = 34394
$ws.Range("J42").Value = 34394
$ws.Range("L42").Value = 34394
$ws.Range("N42").Value = -35364
$ws.Range("H70").Value = 428002.25
$ws.Range("I70").Value = 501000
$ws.Range("K70").Value = 501000
$ws.Range("M70").Value = -500730
$ws.Range("H73").Value = 428002.25
$ws.Range("I73").Value = 501000
$ws.Range("K73").Value = 501000
$ws.Range("M73").Value = -500064
$ws.Range("H80").Value = 55558376
$ws.Range("I80").Value = 111113870
$ws.Range("J80").Value = 2879.2222
$ws.Range("K80").Value = 111113870
$ws.Range("L80").Value = 2879.2222
$ws.Range("M80").Value = -111112872
$ws.Range("N80").Value = -4875.2222
$ws.Range("H83").Value = 55558376
$ws.Range("I83").Value = 111113870
$ws.Range("J83").Value = 2879.2222
$ws.Range("K83").Value = 555569350
$ws.Range("L83").Value = 14396.111
$ws.Range("M83").Value = -555564358
$ws.Range("N83").Value = -24380.111
$ws.Range("H115").Value = 34394
$ws.Range("J115").Value = 34394
$ws.Range("L115").Value = 34394
$ws.Range("N115").Value = -36744
$ws.Range("H122").Value = 20426.555
$ws.Range("I122").Value = 37499.75
$ws.Range("J122").Value = 6768
$ws.Range("K122").Value = 112499.25
$ws.Range("L122").Value = 20304
$ws.Range("M122").Value = -110049.25
$ws.Range("N122").Value = -25204
$ws.Range("H132").Value = 5550
$ws.Range("J132").Value = 5957
$ws.Range("L132").Value = 17871
$ws.Range("N132").Value = -22931

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 37146.145
$ws.Range("I7").Value = 28758.125
$ws.Range("K7").Value = 28758.125
$ws.Range("M7").Value = -28646.125
$ws.Range("H23").Value = 3831.6667
$ws.Range("J23").Value = 3747.5
$ws.Range("L23").Value = 3747.5
$ws.Range("N23").Value = -4207.5
$ws.Range("H30").Value = 208
$ws.Range("I30").Value = 208
$ws.Range("K30").Value = 208
$ws.Range("M30").Value = -100
$ws.Range("H93").Value = 1809.7894
$ws.Range("I93").Value = 1509.2
$ws.Range("K93").Value = 1509.2
$ws.Range("M93").Value = -261.2
$ws.Range("H126").Value = 37146.145
$ws.Range("I126").Value = 28758.125
$ws.Range("K126").Value = 86274.375
$ws.Range("M126").Value = -83804.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 174932.86
$ws.Range("J46").Value = 174932.86
$ws.Range("L46").Value = 174932.86
$ws.Range("N46").Value = -175394.86
$ws.Range("H100").Value = 5495303
$ws.Range("I100").Value = 8929093
$ws.Range("J100").Value = 1238.8
$ws.Range("K100").Value = 17858186
$ws.Range("L100").Value = 2477.6
$ws.Range("M100").Value = -17857645
$ws.Range("N100").Value = -3559.6
$ws.Range("H113").Value = 945.0909
$ws.Range("I113").Value = 776.7143
$ws.Range("K113").Value = 2330.1429
$ws.Range("M113").Value = -160.1428999999998
$ws.Range("H126").Value = 14835.2
$ws.Range("I126").Value = 1419.6666
$ws.Range("K126").Value = 4258.9998
$ws.Range("M126").Value = -1788.9998
$ws.Range("H132").Value = 1319319.4
$ws.Range("I132").Value = 1480.5416
$ws.Range("J132").Value = 4833556.5
$ws.Range("K132").Value = 4441.6248
$ws.Range("L132").Value = 14500669.5
$ws.Range("M132").Value = -1911.6248
$ws.Range("N132").Value = -14505729.5
$ws.Range("H134").Value = 174932.86
$ws.Range("J134").Value = 174932.86
$ws.Range("L134").Value = 524798.58
$ws.Range("N134").Value = -529868.58
$ws.Range("H136").Value = 4015.923
$ws.Range("I136").Value = 1954.2727
$ws.Range("K136").Value = 5862.8181
$ws.Range("M136").Value = -3312.8181
